# Updated symbol list on Thu Dec 22 20:15:57 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price/value cells (column D) — stored as text in the source data,
# so force text format before assigning to avoid Excel's numeric auto-coercion
# (these values rely on exact text formatting, e.g. trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.77"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.64"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.434"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05705"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.420"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.267"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8089"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.069"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1416"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07288"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03089"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03113"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09364"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.917"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001572"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04802"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005806"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006273"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004072"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009911"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.732"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.155"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3263"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1300"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03824"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006651"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1051"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002798"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006497"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005597"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3897"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002098"

# "Hora" column (G) — every row moves from 19 to 20; also numeric-looking text.
$gRange = $ws.Range("G2:G51")
$gRange.NumberFormat = "@"
$gRange.Value = "20"

# Row 8/9 swapped to MXToken/FTXToken with updated link, volume label, and price;
# plus assorted Volume(1h) label tweaks elsewhere (plain text columns B, C, E).
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E27").Value = "26UpBotsUBXTWorstin24h"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E44").Value = "43LocalTradersLCT"
